# Updates cryptos list values/percentages and two row re-orderings (B/C columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.088.45'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '1.892.32'
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("D4").Value = '''0.9972'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.49%  '

$ws.Range("D5").Value = '''0.7432'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '

$ws.Range("D6").Value = '''243.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").Value = '''0.9980'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("D8").Value = '''0.3173'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.80%  '

$ws.Range("D9").Value = '''0.07250'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.91%  '

$ws.Range("D10").Value = '''25.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.29%  '

$ws.Range("D11").Value = '''0.08359'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.04%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '''0.7604'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.917.41'
$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("D14").Value = '''5.427'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.26%  '

$ws.Range("D15").Value = '''92.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.45%  '

$ws.Range("D16").Value = '''6.172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("D17").Value = '30.088.31'
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = '''250.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.93%  '

$ws.Range("D19").Value = '''13.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").Value = '''0.000007870'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''0.9994'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.139.12'
$ws.Range("E22").Value = '  -0.74%  '

$ws.Range("D23").Value = '''8.043'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").Value = '''0.9972'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("D25").Value = '''0.1586'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("D26").Value = '''9.315'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("D27").Value = '''164.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '

$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("D29").Value = '''2.062'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.76%  '

$ws.Range("D30").Value = '''1.478'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").Value = '''4.610'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.18%  '

$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").Value = '''4.239'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.58%  '

$ws.Range("D34").Value = '''0.05383'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = '''1.257'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.79%  '

$ws.Range("D36").Value = '''0.7643'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.99%  '

$ws.Range("D37").Value = '''1.001'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("D38").Value = '''2.722'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").Value = '''0.01975'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.92%  '

$ws.Range("D40").Value = '''2.767'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").Value = '''0.4573'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.81%  '

$ws.Range("D42").Value = '1.102.47'
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''6.090'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''73.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.98%  '

$ws.Range("D45").Value = '''0.8716'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.73%  '

$ws.Range("D46").Value = '''104.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.23%  '

$ws.Range("D47").Value = '''0.9992'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("D48").Value = '''1.872'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").Value = '''7.637'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").Value = '''9.610'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.98%  '

$ws.Range("D51").Value = '2.057.13'
$ws.Range("E51").Value = '  +0.58%  '
